# Generate Report for Handback
# Updates the "Ready for handoff" status rows (for the 27071075-... file) to
# reflect a failed handback transform, records the error detail message on
# the zh-cn and de-de sheets, and widens the "Error Detail" column so the
# message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 corresponds to 27071075-8848-4d23-a425-07606f580f70.md
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: row 3 Status column (C) and Error Detail column (P)
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = "Handback file name: ks31pnah.2jk is different with handoff file name: 27071075-8848-4d23-a425-07606f580f70.3e9abea58123c99540ba4e24e6844394c33fefbd.zh-cn."

# de-de sheet: row 3 Status column (C) and Error Detail column (P)
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = "Handback file name: ks31pnah.2jk is different with handoff file name: 27071075-8848-4d23-a425-07606f580f70.3e9abea58123c99540ba4e24e6844394c33fefbd.de-de."

# Widen the "Error Detail" column (column P, the 16th column) on both
# language sheets so the new error message is fully visible (target width
# of 40 characters). Column A on these sheets is already authored at width
# 40, so read its effective COM ColumnWidth and reuse it to land on exactly
# the same stored width, avoiding COM's internal pixel-rounding drift.
$refWidth = $zhcn.Columns.Item(1).ColumnWidth
$zhcn.Columns.Item(16).ColumnWidth = $refWidth
$dede.Columns.Item(16).ColumnWidth = $refWidth
